$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update A2 with more detail (was "Link without context")
$ws.Range("A2").Value = '"Click Here" Link without context'

# B2 stays "Home" (value unchanged)
$ws.Range("B2").Value = "Home"

# Move the active selection from D2 to A2
$ws.Range("A2").Select()

$wb.Save()
